$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Delete()
$ws.Columns.Item(1).Delete()

$ws.Range("C2").Value = 3
$ws.Range("E2").Value = 6
